$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from H1 (bold, bordered, centered) to the new
# header cells I1:J1, then set their values.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for the new columns I (I0) and J (IF), rows 2-18
$data = @(
    @(7, 8),
    @(9, 9),
    @(9, 9),
    @(6, 7),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(5, 5),
    @(8, 8),
    @(5, 6),
    @(6, 6),
    @(7, 8),
    @(9, 9),
    @(8, 8),
    @(7, 7),
    @(7, 7),
    @(6, 6)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
